# Update "想去人数" (want-to-go count) figures in the "展览" (Exhibition) sheet
# and mirror the same updates in the "全部类型" (All Types) aggregate sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1866
$ws1.Range("F5").Value = 3251
$ws1.Range("F7").Value = 4739
$ws1.Range("F8").Value = 449
$ws1.Range("F10").Value = 158
$ws1.Range("F11").Value = 615
$ws1.Range("F21").Value = 4696
$ws1.Range("F23").Value = 31
$ws1.Range("F25").Value = 5844
$ws1.Range("F27").Value = 1183
$ws1.Range("F29").Value = 653
$ws1.Range("F30").Value = 4409
$ws1.Range("F34").Value = 816
$ws1.Range("F35").Value = 59
$ws1.Range("F36").Value = 740
$ws1.Range("F37").Value = 760

$ws2 = $wb.Worksheets.Item("全部类型")
$ws2.Range("F7").Value = 1866
$ws2.Range("F9").Value = 3251
$ws2.Range("F11").Value = 4739
$ws2.Range("F12").Value = 449
$ws2.Range("F14").Value = 158
$ws2.Range("F15").Value = 615
$ws2.Range("F26").Value = 4696
$ws2.Range("F28").Value = 31
$ws2.Range("F30").Value = 5844
$ws2.Range("F32").Value = 1183
$ws2.Range("F34").Value = 653
$ws2.Range("F35").Value = 4409
$ws2.Range("F40").Value = 816
$ws2.Range("F41").Value = 59
$ws2.Range("F42").Value = 740
$ws2.Range("F43").Value = 760
